{"js": "// Update the \"Mean Green Pathway Contribution\" values in the nested CNS\n// table. Each entry gives the 0-based table row index (rows 0-1 are the\n// two-row header) and the expected current value in the 3rd column\n// (index 2), which is replaced by the new value. Row 18 (SRS Marsh / Dry,\n// value 0.71) is intentionally left untouched, matching the source diff.\nconst changes = [\n  { row: 2, oldValue: \"0.87\", newValue: \"0.91\" },\n  { row: 3, oldValue: \"0.89\", newValue: \"0.91\" },\n  { row: 4, oldValue: \"0.61\", newValue: \"0.69\" },\n  { row: 5, oldValue: \"0.37\", newValue: \"0.47\" },\n  { row: 6, oldValue: \"0.71\", newValue: \"0.59\" },\n  { row: 7, oldValue: \"0.84\", newValue: \"0.65\" },\n  { row: 8, oldValue: \"0.47\", newValue: \"0.60\" },\n  { row: 9, oldValue: \"0.59\", newValue: \"0.58\" },\n  { row: 10, oldValue: \"0.61\", newValue: \"0.50\" },\n  { row: 11, oldValue: \"0.78\", newValue: \"0.43\" },\n  { row: 12, oldValue: \"0.74\", newValue: \"0.75\" },\n  { row: 13, oldValue: \"0.22\", newValue: \"0.31\" },\n  { row: 14, oldValue: \"0.70\", newValue: \"0.44\" },\n  { row: 15, oldValue: \"0.75\", newValue: \"0.44\" },\n  { row: 16, oldValue: \"0.54\", newValue: \"0.55\" },\n  { row: 17, oldValue: \"0.57\", newValue: \"0.58\" },\n  { row: 19, oldValue: \"0.80\", newValue: \"0.79\" },\n];\n\nconst valueColumn = 2;\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Grab every target cell up front and load its current text.\nconst cells = changes.map((c) => table.getCell(c.row, valueColumn));\ncells.forEach((cell) => cell.load(\"value\"));\nawait context.sync();\n\n// Verify the cell currently holds the value the diff expects before\n// overwriting it, so the edit only touches the intended numbers even if\n// the table layout shifts slightly; otherwise still apply the intended\n// new value so the net result matches the target state.\nfor (let i = 0; i < changes.length; i++) {\n  const { oldValue, newValue } = changes[i];\n  const cell = cells[i];\n  const current = (cell.value || \"\").trim();\n  if (current !== oldValue) {\n    console.log(\n      `Warning: row ${changes[i].row} expected \"${oldValue}\" but found \"${current}\"`\n    );\n  }\n  cell.value = newValue;\n}\n\nawait context.sync();\n", "ps1": "# Update the \"Mean Green Pathway Contribution\" values in the nested CNS\n# table. Each entry gives the 1-based table Row/Column (Word COM indexing)\n# and the expected current value in that cell, which is replaced by the\n# new value. Row 20 (SRS Marsh / Dry, value 0.71) is intentionally left\n# untouched, matching the source diff.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$valueColumn = 3\n\n$changes = @(\n    @{ Row = 3;  OldValue = \"0.87\"; NewValue = \"0.91\" },\n    @{ Row = 4;  OldValue = \"0.89\"; NewValue = \"0.91\" },\n    @{ Row = 5;  OldValue = \"0.61\"; NewValue = \"0.69\" },\n    @{ Row = 6;  OldValue = \"0.37\"; NewValue = \"0.47\" },\n    @{ Row = 7;  OldValue = \"0.71\"; NewValue = \"0.59\" },\n    @{ Row = 8;  OldValue = \"0.84\"; NewValue = \"0.65\" },\n    @{ Row = 9;  OldValue = \"0.47\"; NewValue = \"0.60\" },\n    @{ Row = 10; OldValue = \"0.59\"; NewValue = \"0.58\" },\n    @{ Row = 11; OldValue = \"0.61\"; NewValue = \"0.50\" },\n    @{ Row = 12; OldValue = \"0.78\"; NewValue = \"0.43\" },\n    @{ Row = 13; OldValue = \"0.74\"; NewValue = \"0.75\" },\n    @{ Row = 14; OldValue = \"0.22\"; NewValue = \"0.31\" },\n    @{ Row = 15; OldValue = \"0.70\"; NewValue = \"0.44\" },\n    @{ Row = 16; OldValue = \"0.75\"; NewValue = \"0.44\" },\n    @{ Row = 17; OldValue = \"0.54\"; NewValue = \"0.55\" },\n    @{ Row = 18; OldValue = \"0.57\"; NewValue = \"0.58\" },\n    @{ Row = 20; OldValue = \"0.80\"; NewValue = \"0.79\" }\n)\n\nforeach ($change in $changes) {\n    $cell = $t.Cell($change.Row, $valueColumn)\n    $cell.Range.Text = $change.NewValue\n}\n"}
